# teaching-experience_es.xlsx — fix the "2020 - Present" (2020 - Actualmente)
# row and normalize the dash spacing on the other date-range cells in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Universidad El Bosque — Profesora Asociada: was "2020-2024" (wrong, this is
# the current/ongoing position) -> "2020 - Actualmente" ("2020 - Present")
$ws.Range("B2").Value = "2020 - Actualmente"

# Universidad El Bosque — Profesora Asistente: "2015-2019" -> "2015 - 2019"
$ws.Range("B4").Value = "2015 - 2019"

# Universidad de San Buenaventura de Medellín - Sede Ibagué — Profesora
# tiempo completo: "2013-2014" -> "2013 - 2014"
$ws.Range("B9").Value = "2013 - 2014"

# Universidad Antonio Nariño — Profesora catedrática: "2013-2014" -> "2013 - 2014"
$ws.Range("B14").Value = "2013 - 2014"

# Match the author's final cursor position in the saved file.
$ws.Range("C18").Select() | Out-Null
